# Weekly update: insert one new daily price record for "Ajo" (garlic) at
# Feria Lagunitas de Puerto Montt, on 2021-08-30 (Excel serial 44438).
#
# The new record is inserted as row 89, pushing the existing rows 89-140
# down to 90-141 (matching the sheet's existing newest-first-ish ordering
# seen in the diff, where every subsequent row's data equals the row
# that used to sit one position above it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 89..140 down to 90..141 and leave a blank row 89 to fill in.
$ws.Rows.Item(89).Insert()

# Populate the newly inserted row 89 with the new observation.
$ws.Cells.Item(89, 1).Value  = 4
$ws.Cells.Item(89, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(89, 3).Value  = "Los Lagos"
$ws.Cells.Item(89, 4).Value  = 44438
$ws.Cells.Item(89, 5).Value  = 10
$ws.Cells.Item(89, 6).Value  = 100112003
$ws.Cells.Item(89, 7).Value  = "Ajo"
$ws.Cells.Item(89, 8).Value  = "Chino"
$ws.Cells.Item(89, 9).Value  = "Primera"
$ws.Cells.Item(89, 10).Value = 100
$ws.Cells.Item(89, 11).Value = 18000
$ws.Cells.Item(89, 12).Value = 18000
$ws.Cells.Item(89, 13).Value = 18000
$ws.Cells.Item(89, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(89, 15).Value = "China"
$ws.Cells.Item(89, 16).Value = 1800
$ws.Cells.Item(89, 17).Value = 10
$ws.Cells.Item(89, 18).Value = "Hortaliza"

# Make sure the date column keeps the same date-time number format as the
# rest of column D (Insert() already carries this down, but set it
# explicitly so the new row matches regardless of host behaviour).
$ws.Cells.Item(89, 4).NumberFormat = $ws.Cells.Item(90, 4).NumberFormat
